$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.030.93"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "3.366.96"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'570.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'135.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.364.44"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("E11").Value = "  -3.88%  "
$ws.Range("E12").Value = "  -3.06%  "
$ws.Range("D13").Value = "3.941.30"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'25.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("E16").Value = "  -4.72%  "
$ws.Range("D17").Value = "3.365.63"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "61.158.57"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "'13.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'5.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.87%  "
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").Value = "'376.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "'0.554"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("D24").Value = "3.496.17"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").Value = "'70.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'1.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.23%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "'8.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "'2.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'23.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'5.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.26%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'1.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'164.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").Value = "'0.0752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'1.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "'41.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "'23.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.88%  "
$ws.Range("D48").Value = "'23.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  -2.74%  "
$ws.Range("D50").Value = "2.346.74"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "'2.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.76%  "
